$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet2")

# Update the selection/input cells in row 3.
$ws.Range("B3").Value = "All"
$ws.Range("C3").Value = "All"
$ws.Range("E3").Value = "Roll Count"
$ws.Range("F3").Value = "All"

$wb.Save()
